# Apply the 2022FSAdates.xlsx edits described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91: SplitGender 1 -> 0 ---
$ws.Range("E91").Value = 0

# --- Row 92: Category U15 -> U1720, Weapon S -> F ---
$ws.Range("B92").Value = "U1720"
$ws.Range("C92").Value = "F"

# --- Row 93: Category U13T -> U15, Weapon F -> S, Cancelled 0 -> 1,
#             Time 10:30 -> 10:00, Rollcall 10:15 -> 9:45 ---
$ws.Range("B93").Value = "U15"
$ws.Range("C93").Value = "S"
$ws.Range("D93").Value = 1
$ws.Range("F93").Value = "10:00"
$ws.Range("G93").Value = "9:45"

# --- Row 94: SplitGender 1 -> 0 ---
$ws.Range("E94").Value = 0

# --- Row 95: Cancelled 0 -> 1 ---
$ws.Range("D95").Value = 1

# --- Row 96: SplitGender 1 -> 0 ---
$ws.Range("E96").Value = 0

# --- Row 97: Cancelled 0 -> 1 ---
$ws.Range("D97").Value = 1

# --- Row 98: SplitGender 1 -> 0 ---
$ws.Range("E98").Value = 0

# --- Sheet view: active selection ---
$ws.Activate()
$ws.Range("L92").Select()
